$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.752.66"
$ws.Range("E2").Value = "  -0.65%  "

$ws.Range("D3").Value = "3.510.27"
$ws.Range("E3").Value = "  -1.17%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'586.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.96%  "

$ws.Range("D6").Value = "'133.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.14%  "

$ws.Range("D7").Value = "3.510.41"
$ws.Range("E7").Value = "  -1.12%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "'0.489"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.86%  "

$ws.Range("D10").Value = "'0.124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.29%  "

$ws.Range("D11").Value = "'7.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.49%  "

$ws.Range("D12").Value = "'0.385"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.48%  "

$ws.Range("D13").Value = "4.120.75"
$ws.Range("E13").Value = "  -0.82%  "

$ws.Range("D14").Value = "'27.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.63%  "

$ws.Range("D15").Value = "'0.0000180"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.29%  "

$ws.Range("E16").Value = "  +0.62%  "

$ws.Range("D17").Value = "3.523.02"
$ws.Range("E17").Value = "  -0.89%  "

$ws.Range("D18").Value = "64.829.13"
$ws.Range("E18").Value = "  -0.72%  "

$ws.Range("D19").Value = "'9.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.39%  "

$ws.Range("D20").Value = "'14.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.21%  "

$ws.Range("D21").Value = "'5.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.01%  "

$ws.Range("D22").Value = "'389.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.29%  "

$ws.Range("D23").Value = "'0.575"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.45%  "

$ws.Range("D24").Value = "3.660.27"
$ws.Range("E24").Value = "  -0.95%  "

$ws.Range("D25").Value = "'74.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.09%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").Value = "'0.0000109"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.60%  "

$ws.Range("E28").Value = "  +1.45%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "'7.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.53%  "

$ws.Range("D31").Value = "'2.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.91%  "

$ws.Range("D32").Value = "'8.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.16%  "

$ws.Range("D33").Value = "3.522.74"
$ws.Range("E33").Value = "  -0.93%  "

$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").Value = "'23.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.24%  "

$ws.Range("D36").Value = "'0.145"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.65%  "

$ws.Range("D38").Value = "'171.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.29%  "

$ws.Range("D39").Value = "'5.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.18%  "

$ws.Range("D40").Value = "'6.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.32%  "

$ws.Range("D41").Value = "'0.0807"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.25%  "

$ws.Range("D42").Value = "'0.816"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.07%  "

$ws.Range("D43").Value = "'26.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.83%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'42.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.95%  "

$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").Value = "'1.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.50%  "

$ws.Range("D47").Value = "'4.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.89%  "

$ws.Range("D48").Value = "'1.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.36%  "

$ws.Range("D49").Value = "2.491.39"
$ws.Range("E49").Value = "  +1.33%  "

$ws.Range("D50").Value = "'6.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.85%  "

$ws.Range("D51").Value = "'0.896"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.42%  "
